$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price / 1h volume change)
# Values are written with a leading apostrophe to force text, then the
# cell style is reset to "Normal" so no stray NumberFormat/quote-prefix
# style survives on the cell (matches original inlineStr text cells).

# Row 2
$ws.Range("D2").Value = "'60.492.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.50%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.682.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.40%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'519.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'148.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.63%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.89%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.699.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.00%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'6.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.32%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  -0.01%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.343"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.37%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  +1.03%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'3.144.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.00%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'60.448.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.55%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'21.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.72%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "'  +0.25%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.673.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.31%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'353.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.77%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'4.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.17%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'10.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.47%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'6.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.35%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.02%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'63.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.60%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.73%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  +2.85%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.991"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.27%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'0.0₃0828"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.43%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.09%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'6.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.72%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  +0.13%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'19.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.12%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +0.12%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'150.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.45%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'4.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.27%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.954"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -9.70%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'1.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.90%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  +8.70%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.885"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.94%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  +0.88%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'3.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'284.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.61%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.0993"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.41%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'20.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.77%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = "'Mantle"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.611"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.52%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "'FirstDigitalUSD"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.10%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'2.106.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.60%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("B48").Value = "'RenderToken"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'4.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.55%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("B49").Value = "'Hedera"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0539"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.47%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("B50").Value = "'VeChain"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0234"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.05%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'19.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.12%  "
$ws.Range("E51").Style = "Normal"
